# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1406
$wsExpo.Range("F3").Value = 2957
$wsExpo.Range("F4").Value = 24
$wsExpo.Range("F5").Value = 271

# Sheet "全部类型" (all types) - contains the same rows duplicated
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1406
$wsAll.Range("F4").Value = 2957
$wsAll.Range("F5").Value = 24
$wsAll.Range("F7").Value = 271
